$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Delete the entire "draft instructions" block at the top of the document:
#    everything from paragraph 1 up to (and including) the paragraph ending
#    "...repository and submit a link to it." (second copy, with the bolded
#    "Data" run) together with the run of blank paragraphs that trail it.
#    That is paragraphs 1..21 (the paragraph right before the
#    "Battle of Zip codes - Charlotte Edition" title).
# ---------------------------------------------------------------------------
$titlePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Battle of Zip codes " + [char]0x2013 + " Charlotte Edition") {
        $titlePara = $p
        break
    }
}

if ($titlePara -ne $null) {
    $startRng = $d.Paragraphs(1).Range
    $endRng = $titlePara.Range
    $deleteRange = $d.Range($startRng.Start, $endRng.Start)
    $deleteRange.Delete()
}

# ---------------------------------------------------------------------------
# 2. Remove the old hidden "_GoBack" bookmark (it used to sit in the
#    paragraph that begins "The data that was collected...").
# ---------------------------------------------------------------------------
foreach ($bm in $d.Bookmarks) {
    if ($bm.Name -eq "_GoBack") {
        $bm.Delete()
    }
}

# ---------------------------------------------------------------------------
# 3. Re-find the title paragraph (indices shifted after the deletion) and:
#      a) add a fresh "_GoBack" bookmark collapsed at its very start
#      b) strip the stale <w:lastRenderedPageBreak/> cached in its run,
#         while keeping all of the run's character formatting intact.
# ---------------------------------------------------------------------------
$titlePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Battle of Zip codes " + [char]0x2013 + " Charlotte Edition") {
        $titlePara = $p
        break
    }
}

$titleStart = $d.Range($titlePara.Range.Start, $titlePara.Range.Start)
$d.Bookmarks.Add("_GoBack", $titleStart)

$titleTextRng = $d.Range($titlePara.Range.Start, $titlePara.Range.End - 1)
$roundTrip = $titleTextRng.FormattedText
$titleTextRng.FormattedText = $roundTrip

# ---------------------------------------------------------------------------
# 4. Renumber the hyperlink relationship ids that shift down now that the
#    embedded ActiveX control's image/control relationships are gone
#    (rId7->rId5, rId8->rId6, rId9->rId7, rId10->rId8).
# ---------------------------------------------------------------------------
$targets = @(
    "http://data.charlottenc.gov/datasets/cmpd-officer-involved-shootings-individuals-1/geoservice",
    "http://data.charlottenc.gov/datasets/cmpd-officer-involved-shootings-officers-1/geoservice",
    "http://data.charlottenc.gov/datasets/cmpd-officer-involved-shootings-incidents-1/geoservice",
    "http://data.charlottenc.gov/datasets/existing-shopping-centers/geoservice"
)
foreach ($addr in $targets) {
    foreach ($h in $d.Hyperlinks) {
        if ($h.Address -eq $addr) {
            $h.Address = $addr
        }
    }
}

Write-Output "edit complete"
